# Atualização de bases das ligas, do dia: 17-03-2024 às 10:24
#
# The source data rows for a handful of fixtures were re-synced from the
# upstream feed: a few rows had their (id-independent) data swapped /
# rotated with neighbouring rows, and a batch of rows that had been
# scraped prematurely (before the match odds were finalised) were removed
# from the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-sync mismatched rows -------------------------------------------
# Row 95 <-> Row 96: swap everything except the leading id column (A).
$v95 = $ws.Range("B95:AC95").Value2
$v96 = $ws.Range("B96:AC96").Value2
$ws.Range("B95:AC95").Value = $v96
$ws.Range("B96:AC96").Value = $v95

# Row 110 <-> Row 111: swap everything except the leading id column (A).
$v110 = $ws.Range("B110:AC110").Value2
$v111 = $ws.Range("B111:AC111").Value2
$ws.Range("B110:AC110").Value = $v111
$ws.Range("B111:AC111").Value = $v110

# Rows 129, 130, 131: rotate the data up by one (129 <- 130 <- 131 <- 129),
# again leaving the id column (A) untouched.
$v129 = $ws.Range("B129:AC129").Value2
$v130 = $ws.Range("B130:AC130").Value2
$v131 = $ws.Range("B131:AC131").Value2
$ws.Range("B129:AC129").Value = $v130
$ws.Range("B130:AC130").Value = $v131
$ws.Range("B131:AC131").Value = $v129

# --- Drop stale trailing rows -------------------------------------------
# Rows 211-216 were removed from the source feed entirely.
$ws.Rows("211:216").Delete()
